$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Methods")

# Finish "User methods" section: mark rows 71-77 (column A) as "Done"
# using the same "Good" cell style already used on rows 6-70.
$userDoneRows = 71,72,73,74,75,76,77
foreach ($r in $userDoneRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "Done"
    $cell.Style = "Good"
}

# Add "Network methods" section completion: mark rows 81-84 (column A)
# as "Done" as well, leaving rows 85-88 (Inbox/Sites) unmarked.
$networkDoneRows = 81,82,83,84
foreach ($r in $networkDoneRows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = "Done"
    $cell.Style = "Good"
}

# Update the current selection/active cell to reflect where work left off.
$ws.Activate()
$ws.Cells.Item(85, 1).Select()
